{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n// the \"\u00a9 2020 . Contact: ...\" footer paragraph, and the blank paragraph that\n// separates them from the preceding \"LOQ4233: ...\" requirement line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOQ4233: ...\") and the two text paragraphs\n// that must be removed, by content rather than a hard-coded index so the\n// script is resilient to unrelated edits elsewhere in the document.\nlet anchorIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"LOQ4233\") !== -1) {\n    anchorIndex = i;\n  }\n  if (text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex !== -1 && copyrightIndex !== -1) {\n  // Delete the blank separator paragraph right before \"Ver no Jupiter ...\"\n  // (only if it really sits between the anchor and that paragraph and is\n  // blank), then delete the two content paragraphs themselves.\n  if (\n    anchorIndex !== -1 &&\n    jupiterIndex === anchorIndex + 2 &&\n    items[anchorIndex + 1].text.trim() === \"\"\n  ) {\n    items[anchorIndex + 1].delete();\n  } else if (jupiterIndex > 0 && items[jupiterIndex - 1].text.trim() === \"\") {\n    items[jupiterIndex - 1].delete();\n  }\n\n  items[jupiterIndex].delete();\n  items[copyrightIndex].delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph,\n# the \"\u00a9 2020 . Contact: ... Powered by Jekyll ...\" footer paragraph, and the\n# blank paragraph that separates them from the preceding\n# \"LOQ4233: ...\" requirement line.\n$d = $word.ActiveDocument\n\n# Locate the two content paragraphs to remove by searching their text rather\n# than relying on a fixed paragraph index.\n$jupiterRange = $d.Content\n$jupiterFound = $jupiterRange.Find.Execute(\"Ver no Jupiter\")\n$idxJupiter = $jupiterRange.Paragraphs.Item(1).Index\n\n$copyrightRange = $d.Content\n$copyrightFound = $copyrightRange.Find.Execute(\"Powered by Jekyll\")\n$idxCopyright = $copyrightRange.Paragraphs.Item(1).Index\n\nif ($jupiterFound -and $copyrightFound) {\n    $jupiterPara = $d.Paragraphs.Item($idxJupiter)\n    $prevPara = $jupiterPara.Previous()\n    $prevIdx = $prevPara.Index\n    $prevIsBlank = ($prevPara.Range.Text.Trim() -eq \"\")\n\n    # Delete highest index first so the lower indices stay valid.\n    $d.Paragraphs.Item($idxCopyright).Range.Delete()\n    $d.Paragraphs.Item($idxJupiter).Range.Delete()\n    if ($prevIsBlank) {\n        $d.Paragraphs.Item($prevIdx).Range.Delete()\n    }\n}\n"}
